$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = "Fixed ROS with Realsense SDK, started ROS research"
$ws.Range("C24").Value = "Finished getTF, updated class diagram"

$ws.Range("D24").Select()
